$wb = $excel.ActiveWorkbook

# --- Sheet: 土地 ---
$ws = $wb.Worksheets.Item("土地")
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "area"
$ws.Range("D1").Value = "share_portion"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"
$ws.Range("H1").Value = "acquire_value"
$ws.Range("I1").Value = "property_category"
$ws.Range("J1").Value = "category"
$ws.Range("K1").Value = "date"
$ws.Range("L1").Value = "legislator_name"
$ws.Range("M1").Value = "legislator_id"
$ws.Range("N1").Value = "source_file"
$ws.Range("O1").Value = "index"
$ws.Range("B2").Value = "桃圜縣中壢市石頭段00490032地號"
$ws.Range("F2").Value = "80年08月07日"
$ws.Range("I2").Value = "land"
$ws.Range("J2").Value = "normal"
$ws.Range("K2").Value = "2013-12-20"
$ws.Range("L2").Value = "廖正井"
$ws.Range("M2").Value = 1711
$ws.Range("N2").Value = "tmp393a1"
$ws.Range("O2").Value = 13
$ws.Range("B3").Value = "桃圜縣中壢市石頭段00490033地號"
$ws.Range("F3").Value = "80年08月07日"
$ws.Range("I3").Value = "land"
$ws.Range("J3").Value = "normal"
$ws.Range("K3").Value = "2013-12-20"
$ws.Range("L3").Value = "廖正井"
$ws.Range("M3").Value = 1711
$ws.Range("N3").Value = "tmp393a1"
$ws.Range("O3").Value = 14
$ws.Range("B4").Value = "桃圜縣中壢市石頭段01320055地號"
$ws.Range("F4").Value = "80年08月07日"
$ws.Range("I4").Value = "land"
$ws.Range("J4").Value = "normal"
$ws.Range("K4").Value = "2013-12-20"
$ws.Range("L4").Value = "廖正井"
$ws.Range("M4").Value = 1711
$ws.Range("N4").Value = "tmp393a1"
$ws.Range("O4").Value = 15
$ws.Range("B5").Value = "臺北市大安區大安段二小段01110000地號"
$ws.Range("D5").Value = "10000分之566"
$ws.Range("F5").Value = "81年08月12日"
$ws.Range("I5").Value = "land"
$ws.Range("J5").Value = "normal"
$ws.Range("K5").Value = "2013-12-20"
$ws.Range("L5").Value = "廖正井"
$ws.Range("M5").Value = 1711
$ws.Range("N5").Value = "tmp393a1"
$ws.Range("O5").Value = 16

# --- Sheet: 建物 ---
$ws = $wb.Worksheets.Item("建物")
$ws.Range("B2").Value = "桃圜縣中壢市石頭段06151000建號"
$ws.Range("F2").Value = "80年02月13曰"
$ws.Range("B3").Value = "桃圜縣中壢市石頭段06152000建號"
$ws.Range("F3").Value = "80年02月13曰"
$ws.Range("B4").Value = "臺北市大安區大安段二小段01874000建號"
$ws.Range("F4").Value = "81年10月13日"
$ws.Range("B5").Value = "臺北市大安區大安段二小段01875000建號"
$ws.Range("F5").Value = "81年10月13曰"
$ws.Range("H5").Value = "(超過五年停車位）"

# --- Sheet: 存款 ---
$ws = $wb.Worksheets.Item("存款")
$ws.Range("B5").Value = "台北富邦商業銀行市府分行"
$ws.Range("B6").Value = "台北富邦商業銀行城中分行"
$ws.Range("B7").Value = "中華郵政股份有限公司桃圜府前郵局"
$ws.Range("B8").Value = "中華郵政股份有限公司台北信維郵局"
$ws.Range("B9").Value = "中華郵政股份有限公司台北信維郵局"
$ws.Range("B10").Value = "國泰世華商業銀行信義分行"
$ws.Range("B11").Value = "國泰世華商業銀行信義分行"
$origStyle_F11 = $ws.Range("F11").Style
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "1041062"
$ws.Range("F11").Style = $origStyle_F11
$ws.Range("B12").Value = "國泰世華商業銀行信義分行"
$ws.Range("B16").Value = "中國信託商業銀行敦南分行"
$ws.Range("B17").Value = "台北富邦商業銀行敦和分行"

# --- Sheet: 事業投資 ---
$ws = $wb.Worksheets.Item("事業投資")
$ws.Range("F2").Value = "100年07月14日"
$ws.Range("F3").Value = "101年07月15日"
